$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting so values
# like "23.60" or "0.08610" are not silently coerced into numbers and lose
# trailing zeros / dot-grouping formatting.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "24.135.03"
$ws.Range("E2").Value = "  -2.57%  "

# Row 3
$ws.Range("D3").Value = "1.642.86"
$ws.Range("E3").Value = "  -2.47%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$ws.Range("D5").Value = "307.73"
$ws.Range("E5").Value = "  -1.84%  "

# Row 6
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("E7").Value = "  -0.82%  "

# Row 8
$ws.Range("D8").Value = "0.3858"
$ws.Range("E8").Value = "  -2.91%  "

# Row 9
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "49.47"
$ws.Range("E10").Value = "  -4.40%  "

# Row 11
$ws.Range("D11").Value = "1.354"
$ws.Range("E11").Value = "  -5.07%  "

# Row 12
$ws.Range("D12").Value = "0.08610"
$ws.Range("E12").Value = "  -0.80%  "

# Row 13
$ws.Range("D13").Value = "23.60"
$ws.Range("E13").Value = "  -6.36%  "

# Row 14
$ws.Range("D14").Value = "7.095"
$ws.Range("E14").Value = "  -2.97%  "

# Row 15
$ws.Range("D15").Value = "0.00001288"
$ws.Range("E15").Value = "  -2.48%  "

# Row 16
$ws.Range("D16").Value = "7.455"
$ws.Range("E16").Value = "  -4.41%  "

# Row 17
$ws.Range("D17").Value = "1.631.03"
$ws.Range("E17").Value = "  +2.01%  "

# Row 18
$ws.Range("D18").Value = "94.79"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19
$ws.Range("D19").Value = "0.06901"
$ws.Range("E19").Value = "  -2.66%  "

# Row 20
$ws.Range("D20").Value = "20.37"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21
$ws.Range("D21").Value = "6.896"
$ws.Range("E21").Value = "  -3.40%  "

# Row 22
$ws.Range("E22").Value = "  -0.30%  "

# Row 23
$ws.Range("D23").Value = "13.55"
$ws.Range("E23").Value = "  -3.93%  "

# Row 24
$ws.Range("D24").Value = "24.146.06"

# Row 25
$ws.Range("D25").Value = "2.424"
$ws.Range("E25").Value = "  +2.47%  "

# Row 26
$ws.Range("D26").Value = "2.814"
$ws.Range("E26").Value = "  +1.40%  "

# Row 27
$ws.Range("D27").Value = "22.36"
$ws.Range("E27").Value = "  -6.52%  "

# Row 28
$ws.Range("D28").Value = "157.74"
$ws.Range("E28").Value = "  -2.88%  "

# Row 29
$ws.Range("D29").Value = "8.574"
$ws.Range("E29").Value = "  +9.52%  "

# Row 30
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "5.366"
$ws.Range("E30").Value = "  -6.58%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "140.23"
$ws.Range("E31").Value = "  -6.70%  "

# Row 32
$ws.Range("D32").Value = "2.411"
$ws.Range("E32").Value = "  -7.16%  "

# Row 33
$ws.Range("D33").Value = "1.824.05"
$ws.Range("E33").Value = "  -2.86%  "

# Row 34
$ws.Range("D34").Value = "6.993"
$ws.Range("E34").Value = "  +0.53%  "

# Row 35
$ws.Range("D35").Value = "0.08085"
$ws.Range("E35").Value = "  -4.58%  "

# Row 36
$ws.Range("D36").Value = "0.02896"
$ws.Range("E36").Value = "  -6.00%  "

# Row 37
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2677"
$ws.Range("E37").Value = "  -4.68%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.9478"
$ws.Range("E38").Value = "  -6.60%  "

# Row 39
$ws.Range("D39").Value = "0.09193"

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.457"
$ws.Range("E40").Value = "  -0.87%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "9.972"
$ws.Range("E41").Value = "  -4.83%  "

# Row 42
$ws.Range("D42").Value = "0.7529"
$ws.Range("E42").Value = "  -5.47%  "

# Row 43
$ws.Range("D43").Value = "12.99"
$ws.Range("E43").Value = "  -5.32%  "

# Row 44
$ws.Range("D44").Value = "16.03"
$ws.Range("E44").Value = "  -4.39%  "

# Row 45
$ws.Range("D45").Value = "0.6893"
$ws.Range("E45").Value = "  -3.92%  "

# Row 46
$ws.Range("D46").Value = "2.453"
$ws.Range("E46").Value = "  -5.18%  "

# Row 47
$ws.Range("E47").Value = "  -2.49%  "

# Row 48
$ws.Range("D48").Value = "0.9999"
$ws.Range("E48").Value = "  -0.26%  "

# Row 49
$ws.Range("D49").Value = "0.08375"
$ws.Range("E49").Value = "  -4.50%  "

# Row 50
$ws.Range("D50").Value = "1.258"
$ws.Range("E50").Value = "  -6.34%  "

# Row 51
$ws.Range("D51").Value = "133.11"
$ws.Range("E51").Value = "  -3.83%  "
